# Fix: clear the old sample rows/hyperlink and load the real meeting schedule,
# so the "join" tab can key off day/user without a dead Google Meet hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old hyperlink on C2 (the stale Google Meet URL) before rewriting values.
$ws.Hyperlinks.Delete()

# Drop the stray formatted-but-empty row 9 left over from the old layout.
$ws.Rows("9:9").Clear()

# --- Header row (unchanged labels, but re-apply cleanly) ---
$ws.Range("A1").Value = "day"
$ws.Range("B1").Value = "time"
$ws.Range("C1").Value = "meeting_id"
$ws.Range("D1").Value = "passcode"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "user"

# --- Row 2: monday ---
$ws.Range("A2").Value = "monday"
$ws.Range("B2").Value = 0.66666666666666663
$ws.Range("C2").Value = "943 835 6448"
$ws.Range("D2").Value = 123456
$ws.Range("E2").Value = "zoom"
$ws.Range("F2").Value = "Brenda Carranco"

# --- Row 3: tuesday ---
$ws.Range("A3").Value = "tuesday"
$ws.Range("B3").Value = 0.66666666666666663
$ws.Range("C3").Value = "857 069 3962"
$ws.Range("D3").Value = 536994
$ws.Range("E3").Value = "zoom"
$ws.Range("F3").Value = "Jesus Rodriguez"

# --- Row 4: wednesday ---
$ws.Range("A4").Value = "wednesday"
$ws.Range("B4").Value = 0.54166666666666663
$ws.Range("C4").Value = "650 460 3098"
$ws.Range("D4").Value = 123456
$ws.Range("E4").Value = "zoom"
$ws.Range("F4").Value = "David Haro"

# --- Row 5: saturday ---
$ws.Range("A5").Value = "saturday"
$ws.Range("B5").Value = 0.29166666666666669
$ws.Range("C5").Value = "650 460 3098"
$ws.Range("D5").Value = 123456
$ws.Range("E5").Value = "zoom"
$ws.Range("F5").Value = "David Haro"

# --- Row 6: thursday ---
$ws.Range("A6").Value = "thursday"
$ws.Range("B6").Value = 0.79166666666666663
$ws.Range("C6").Value = 8484246662#
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "055902"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").Value = "zoom"
$ws.Range("F6").Value = "Miguel Alvarez"

# --- Row 7: friday ---
$ws.Range("A7").Value = "friday"
$ws.Range("B7").Value = 0.33333333333333331
$ws.Range("C7").Value = 7544068640#
$ws.Range("D7").Value = 652801
$ws.Range("E7").Value = "zoom"
$ws.Range("F7").Value = "All Hands Meeting"

# --- Formatting pass ---

# Baseline for the whole table: centered, no wrap, regular (non-bold, non-underline) font.
# Also reset font color back to the normal theme (clears the lingering hyperlink-blue
# look that the old style carried on C2/C4/C6 even before the hyperlink itself is gone).
$whole = $ws.Range("A1:F7")
$whole.HorizontalAlignment = -4108
$whole.VerticalAlignment = -4107
$whole.WrapText = $false
$whole.Font.Bold = $false
$whole.Font.Underline = $false
$whole.Font.ThemeColor = 1

# Bold header row, wrapped on the first four columns, plain (no wrap) on the last two.
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").WrapText = $true
$ws.Range("E1:F1").Font.Bold = $true

# Day / meeting_id / passcode columns for the "zoom" rows wrap their text.
$ws.Range("A2:A5").WrapText = $true
$ws.Range("C2:C5").WrapText = $true
$ws.Range("D2:D5").WrapText = $true

# Time column: h:mm, wrapped for the first four data rows, unwrapped for the last two.
$ws.Range("B2:B5").NumberFormat = "h:mm"
$ws.Range("B2:B5").WrapText = $true
$ws.Range("B6:B7").NumberFormat = "h:mm"

# Passcode stored as text (leading zeros) for the thursday row.
$ws.Range("D6").NumberFormat = "@"

# Friday row keeps the underline look (previously the hyperlink font) on day + passcode.
$ws.Range("A7").Font.Underline = $true
$ws.Range("D7").Font.Underline = $true
